$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (Price, Volume(1h)) new values.
# Price values ("D") are forced to text via a leading apostrophe so the
# engine keeps them as plain text cells (matching the source data, which
# stores these as inline/shared strings, not numbers) instead of silently
# re-typing them as numeric cells.
$updates = @{
    2  = @{ D = "76.621.22"; E = "  +1.29%  " }
    3  = @{ D = "2.890.85";  E = "  +8.00%  " }
    4  = @{            E = "  -0.10%  " }
    5  = @{ D = "196.66";    E = "  +4.78%  " }
    6  = @{ D = "600.38";    E = "  +2.21%  " }
    7  = @{            E = "  -0.01%  " }
    8  = @{            E = "  +3.09%  " }
    9  = @{ D = "0.193";     E = "  -0.91%  " }
    10 = @{ D = "2.888.18";  E = "  +7.94%  " }
    11 = @{ D = "0.397";     E = "  +10.75%  " }
    12 = @{            E = "  -1.77%  " }
    13 = @{ D = "4.92";      E = "  +4.24%  " }
    14 = @{ D = "3.391.00";  E = "  +6.97%  " }
    15 = @{ D = "76.476.45"; E = "  +1.39%  " }
    16 = @{ D = "27.50";     E = "  +3.68%  " }
    17 = @{            E = "  +0.65%  " }
    18 = @{ D = "2.890.37";  E = "  +7.93%  " }
    19 = @{            E = "  -2.77%  " }
    20 = @{            E = "  +5.18%  " }
    21 = @{ D = "383.18";    E = "  +2.83%  " }
    22 = @{            E = "  +1.43%  " }
    23 = @{ D = "4.15";      E = "  +1.41%  " }
    24 = @{ D = "71.90";     E = "  +2.64%  " }
    25 = @{ D = "0.998";     E = "  -0.09%  " }
    26 = @{ D = "3.040.09";  E = "  +7.80%  " }
    27 = @{            E = "  +0.88%  " }
    28 = @{ D = "9.78";      E = "  +4.42%  " }
    29 = @{            E = "  +11.51%  " }
    30 = @{ D = "0.999";     E = "  -0.42%  " }
    31 = @{            E = "  +0.66%  " }
    32 = @{ D = "512.90";    E = "  -0.72%  " }
    33 = @{ D = "7.77";      E = "  +0.44%  " }
    34 = @{            E = "  +3.46%  " }
    35 = @{ D = "0.999";     E = "  -0.04%  " }
    36 = @{ D = "168.40";    E = "  +3.01%  " }
    37 = @{ D = "20.14";     E = "  +4.78%  " }
    38 = @{ D = "0.117";     E = "  -1.63%  " }
    39 = @{            E = "  +0.74%  " }
    40 = @{ D = "183.61";    E = "  +8.69%  " }
    41 = @{            E = "  -0.06%  " }
    42 = @{ D = "0.347";     E = "  +4.89%  " }
    43 = @{            E = "  +1.56%  " }
    44 = @{ D = "1.69";      E = "  -0.70%  " }
    45 = @{ D = "0.0925";    E = "  +9.66%  " }
    46 = @{ D = "1.24";      E = "  +3.96%  " }
    47 = @{ D = "40.20";     E = "  +2.38%  " }
    48 = @{ D = "2.37";      E = "  +0.12%  " }
    49 = @{ D = "0.697" }
    50 = @{ D = "0.580";     E = "  +8.42%  " }
    51 = @{            E = "  +3.37%  " }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    if ($vals.ContainsKey("D")) {
        $ws.Range("D$rowNum").Value = "'" + $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $vals["E"]
    }
}
